$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.041.29'
$ws.Range("E2").Value = '  -0.43%  '
$ws.Range("D3").Value = '2.220.14'
$ws.Range("E3").Value = '  -1.32%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.98'
$ws.Range("E5").Value = '  -1.79%  '
$ws.Range("E6").Value = '  +1.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '73.77'
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("E8").Value = '  +0.21%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.617'
$ws.Range("E9").Value = '  -0.36%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '43.73'
$ws.Range("E10").Value = '  +5.92%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0959'
$ws.Range("E11").Value = '  +2.62%  '
$ws.Range("E12").Value = '  +0.28%  '
$ws.Range("E13").Value = '  +0.30%  '
$ws.Range("D14").Value = '2.552.36'
$ws.Range("E14").Value = '  -1.28%  '
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.845'
$ws.Range("E15").Value = '  -0.72%  '
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.25'
$ws.Range("E16").Value = '  -1.25%  '
$ws.Range("D17").Value = '2.232.46'
$ws.Range("E17").Value = '  -0.65%  '
$ws.Range("D18").Value = '41.900.01'
$ws.Range("E18").Value = '  -0.47%  '
$ws.Range("E19").Value = '  +13.15%  '
$ws.Range("E20").Value = '  +0.98%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.41'
$ws.Range("E21").Value = '  +0.80%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.45'
$ws.Range("E22").Value = '  +31.05%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '229.76'
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("E24").Value = '  -6.80%  '
$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.54'
$ws.Range("E25").Value = '  +3.65%  '
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.08%  '
$ws.Range("E27").Value = '  +1.31%  '
$ws.Range("E28").Value = '  -1.48%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.16'
$ws.Range("E29").Value = '  -3.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '166.64'
$ws.Range("E30").Value = '  -1.38%  '
$ws.Range("E31").Value = '  -0.45%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.68'
$ws.Range("E32").Value = '  +16.55%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0798'
$ws.Range("E33").Value = '  -2.98%  '
$ws.Range("E34").Value = '  -0.22%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '29.26'
$ws.Range("E35").Value = '  -3.39%  '
$ws.Range("E36").Value = '  -4.42%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.28'
$ws.Range("E37").Value = '  -4.89%  '
$ws.Range("E38").Value = '  +0.27%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '12.99'
$ws.Range("E39").Value = '  -3.96%  '
$ws.Range("E40").Value = '  -1.74%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '65.34'
$ws.Range("E41").Value = '  +5.57%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.67'
$ws.Range("E42").Value = '  -1.83%  '
$ws.Range("E43").Value = '  -1.75%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.72'
$ws.Range("E44").Value = '  +0.75%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '104.50'
$ws.Range("E45").Value = '  -3.37%  '
$ws.Range("E46").Value = '  +0.59%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.41'
$ws.Range("E47").Value = '  +5.67%  '
$ws.Range("E48").Value = '  +0.44%  '
$ws.Range("E49").Value = '  -0.25%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.71'
$ws.Range("E50").Value = '  +0.67%  '
$ws.Range("D51").Value = '2.426.85'
$ws.Range("E51").Value = '  -1.30%  '
